$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.799.34"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.774.24"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "355.33"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "109.55"
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").Value = "40.04"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("D11").Value = "0.0850"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "19.36"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").Value = "3.212.48"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "2.801.12"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "0.929"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").Value = "51.813.75"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "3.15"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  -4.08%  "
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "274.04"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "69.84"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("D26").Value = "26.46"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").Value = "0.0459"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "33.86"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "0.0843"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").Value = "5.22"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "17.99"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("D41").Value = "2.52"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").Value = "120.98"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("D45").Value = "22.13"
$ws.Range("E45").Value = "  -6.61%  "
$ws.Range("D46").Value = "2.065.15"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("E47").Value = "  -4.31%  "
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("D49").Value = "5.66"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").Value = "0.922"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").Value = "8.93"
$ws.Range("E51").Value = "  -0.02%  "
